# Replace the task list with the new, red-coloured set of items and drop the
# old, now-finished items, while keeping the trailing empty Heading1
# paragraph that carries the _GoBack bookmark.
#
# The whole body is rewritten in one shot via Range.InsertXML so that run
# splits (needed for the "javascript" spell-check markers) and the new
# w:color direct formatting land exactly as Word itself would produce them.

$d = $word.ActiveDocument

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve">Upload markers to </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>javascript</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve"> file</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Make new web page for blog (follow previous google maps one)</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Add markers</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Add menu to zoom to and only turn on markers from the chosen movie</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Add clustering (too many markers web page)</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Add year for each movie</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Add address for each movie</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$d.Content.InsertXML($xml)
